$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DBD")

# Row 22: new field "ActualFilingDate"
$ws.Range("A22").Value = 14
$ws.Range("B22").Value = "ActualFilingDate"
$ws.Range("C22").Value = "實際報送日期"
$ws.Range("D22").Value = "Decimald"
$ws.Range("E22").Value = 8

# Row 23: new field "ActualFilingMark"
$ws.Range("A23").Value = 15
$ws.Range("B23").Value = "ActualFilingMark"
$ws.Range("C23").Value = "實際報送記號"
$ws.Range("D23").Value = "VARCHAR2"
$ws.Range("E23").Value = 3

# Reflect the updated view state (selection moved to A24, scrolled so row 13 is at top)
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A24").Select()
